$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.963.14'
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").Value = '3.527.26'
$ws.Range("E3").Value = '  -1.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.85'
$ws.Range("E5").Value = '  -1.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.87'
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("D7").Value = '3.527.05'
$ws.Range("E7").Value = '  -1.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  -1.70%  '
$ws.Range("E11").Value = '  -2.87%  '
$ws.Range("E12").Value = '  -2.30%  '
$ws.Range("D13").Value = '4.126.49'
$ws.Range("E13").Value = '  -1.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000201'
$ws.Range("E14").Value = '  -3.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.68'
$ws.Range("E15").Value = '  -4.80%  '
$ws.Range("D16").Value = '3.522.32'
$ws.Range("E16").Value = '  -1.82%  '
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").Value = '65.929.56'
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.96'
$ws.Range("E19").Value = '  -4.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.19'
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.39'
$ws.Range("E21").Value = '  -3.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '415.03'
$ws.Range("E22").Value = '  -4.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.598'
$ws.Range("E23").Value = '  -3.40%  '
$ws.Range("E24").Value = '  -2.91%  '
$ws.Range("D25").Value = '3.668.86'
$ws.Range("E25").Value = '  -1.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("E27").Value = '  -3.25%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.84'
$ws.Range("E28").Value = '  -2.78%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.97'
$ws.Range("E29").Value = '  -2.66%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.44'
$ws.Range("E30").Value = '  -2.71%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").Value = '3.525.02'
$ws.Range("E32").Value = '  -1.01%  '
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.39'
$ws.Range("E34").Value = '  -3.95%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.53'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.28'
$ws.Range("E37").Value = '  -12.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '175.91'
$ws.Range("E38").Value = '  +1.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.27'
$ws.Range("E39").Value = '  -6.43%  '
$ws.Range("E40").Value = '  -8.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0820'
$ws.Range("E41").Value = '  -3.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.07'
$ws.Range("E42").Value = '  -2.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.860'
$ws.Range("E43").Value = '  -3.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.31'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.80'
$ws.Range("E45").Value = '  -8.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("E47").Value = '  -4.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.10'
$ws.Range("E48").Value = '  -1.07%  '
$ws.Range("E49").Value = '  -7.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.67'
$ws.Range("E50").Value = '  -3.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.01'
$ws.Range("E51").Value = '  -8.59%  '
